$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns keep their original text formatting
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '30.092.61'
$ws.Range("E2").Value = '  -1.96%  '
$ws.Range("D3").Value = '1.835.41'
$ws.Range("E3").Value = '  -3.25%  '
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.28%  '
$ws.Range("D5").Value = '231.36'
$ws.Range("E5").Value = '  -3.25%  '
$ws.Range("D6").Value = '1.005'
$ws.Range("E6").Value = '  +0.41%  '
$ws.Range("D7").Value = '0.4664'
$ws.Range("E7").Value = '  -3.47%  '
$ws.Range("D8").Value = '0.2687'
$ws.Range("E8").Value = '  -5.74%  '
$ws.Range("D9").Value = '0.06226'
$ws.Range("E9").Value = '  -4.98%  '
$ws.Range("D10").Value = '1.837.19'
$ws.Range("E10").Value = '  -5.75%  '
$ws.Range("D11").Value = '0.07398'
$ws.Range("E11").Value = '  -0.81%  '
$ws.Range("D12").Value = '16.00'
$ws.Range("E12").Value = '  -4.40%  '
$ws.Range("D13").Value = '4.870'
$ws.Range("E13").Value = '  -4.62%  '
$ws.Range("D14").Value = '83.68'
$ws.Range("E14").Value = '  -5.00%  '
$ws.Range("D15").Value = '0.6197'
$ws.Range("E15").Value = '  -7.13%  '
$ws.Range("D16").Value = '30.078.61'
$ws.Range("E16").Value = '  -1.96%  '
$ws.Range("D17").Value = '0.9994'
$ws.Range("E17").Value = '  -0.06%  '
$ws.Range("D18").Value = '12.46'
$ws.Range("E18").Value = '  -6.38%  '
$ws.Range("D19").Value = '224.22'
$ws.Range("E19").Value = '  -3.11%  '
$ws.Range("D20").Value = '0.000007232'
$ws.Range("E20").Value = '  -5.03%  '
$ws.Range("B21").Value = 'BinanceUSD'
$ws.Range("C21").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D21").Value = '1.003'
$ws.Range("E21").Value = '  +0.18%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '4.847'
$ws.Range("E22").Value = '  -8.29%  '
$ws.Range("B23").Value = 'Chainlink'
$ws.Range("C23").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D23").Value = '5.848'
$ws.Range("E23").Value = '  -6.19%  '
$ws.Range("B24").Value = 'Monero'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D24").Value = '166.90'
$ws.Range("E24").Value = '  -1.78%  '
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").Value = '9.127'
$ws.Range("E25").Value = '  -2.24%  '
$ws.Range("B26").Value = 'EthereumClassic'
$ws.Range("C26").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D26").Value = '17.62'
$ws.Range("E26").Value = '  -6.10%  '
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").Value = '1.847'
$ws.Range("E27").Value = '  -6.12%  '
$ws.Range("B28").Value = 'Stellar'
$ws.Range("C28").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D28").Value = '0.1015'
$ws.Range("E28").Value = '  -0.54%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '1.380'
$ws.Range("E29").Value = '  -1.81%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").Value = '4.064'
$ws.Range("E30").Value = '  -6.72%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = '3.783'
$ws.Range("E31").Value = '  -6.16%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = '0.04843'
$ws.Range("E32").Value = '  -5.40%  '
$ws.Range("B33").Value = 'ARBITRUM'
$ws.Range("C33").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D33").Value = '1.139'
$ws.Range("E33").Value = '  -6.48%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = '0.6998'
$ws.Range("E34").Value = '  -7.68%  '
$ws.Range("B35").Value = 'Frax'
$ws.Range("C35").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D35").Value = '0.9993'
$ws.Range("E35").Value = '  -0.37%  '
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").Value = '2.703'
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '0.01824'
$ws.Range("E37").Value = '  -3.31%  '
$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").Value = '2.610'
$ws.Range("E38").Value = '  -1.84%  '
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").Value = '0.9006'
$ws.Range("E39").Value = '  -2.30%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").Value = '1.934'
$ws.Range("E40").Value = '  -7.01%  '
$ws.Range("B41").Value = 'Quant'
$ws.Range("C41").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D41").Value = '104.05'
$ws.Range("E41").Value = '  -2.83%  '
$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").Value = '0.9984'
$ws.Range("E42").Value = '  -0.63%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '5.520'
$ws.Range("E43").Value = '  -3.79%  '
$ws.Range("D44").Value = '0.4042'
$ws.Range("E44").Value = '  -6.06%  '
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").Value = '6.930'
$ws.Range("E45").Value = '  -6.74%  '
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").Value = '0.1180'
$ws.Range("E46").Value = '  -7.56%  '
$ws.Range("D47").Value = '59.05'
$ws.Range("E47").Value = '  -8.69%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '8.453'
$ws.Range("E48").Value = '  -6.06%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '32.86'
$ws.Range("E49").Value = '  -3.03%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '0.05554'
$ws.Range("E50").Value = '  -2.26%  '
$ws.Range("D51").Value = '1.373'
$ws.Range("E51").Value = '  -7.90%  '
